# Apply the edits described by the diff: rewrite the title, author,
# contact line, and all body/summary prose from the biodiversity essay
# to the multicultural-literature essay; then append a trailing empty
# paragraph before the section break, matching the target document.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $null = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- Title ---
Replace-Text ('Unveiling the Diversity of Life: Exploring Biodiversity') ('Through the Lens of Diversity: Navigating the Complexities of Multicultural Literature')

# --- Author name ("Dr. Isabella Medina" -> "Dr. Olivia Bennett") ---
Replace-Text 'Isabella Medina' 'Olivia Bennett'

# --- Contact / email line ---
Replace-Text 'isabella' 'oliviabennett@valid'
Replace-Text 'medina@biodiversity.ac.uk' 'edu'

# --- Main body paragraph (keeps the two manual line breaks) ---
$bodyOld = 'In the vast expanse of our planet, life exists in an astonishing array of forms, from the microscopic wonders of bacteria to the towering giants of the forest. This remarkable diversity, known as biodiversity, encompasses the multitude of species that inhabit Earth''s diverse ecosystems, playing crucial roles in maintaining ecological balance and supporting human well-being. Exploring biodiversity involves unraveling the intricate web of relationships among species and their environments, appreciating the delicate equilibrium that sustains life on Earth.' + [char]11 + '' + [char]11 + 'Biodiversity is a symphony of life, a chorus of interconnected species, each with its unique voice adding to the richness of the composition. From the rainforests of the Amazon to the vibrant coral reefs of the ocean depths, biodiversity manifests in a kaleidoscope of colors, shapes, and behaviors. Each organism, whether a tiny insect or a majestic whale, holds a piece of the puzzle that completes the tapestry of life. Understanding and preserving biodiversity is not merely an academic pursuit; it is a responsibility we owe to future generations, ensuring the continued flourishing of life on Earth.' + [char]11 + '' + [char]11 + 'From its role in maintaining ecosystem stability, providing food, and supporting livelihoods, to its immense cultural and aesthetic value, biodiversity touches every aspect of human existence. Yet, human activities, such as deforestation, pollution, and climate change, are threatening this intricate web of life. By recognizing the intrinsic value of biodiversity and taking collective action to protect and restore ecosystems, we can safeguard the future of life on Earth and ensure the well-being of generations to come.'
$bodyNew = 'Embarking on a literary journey through the rich landscapes of multicultural literature invites us into a realm of kaleidoscopic experiences waiting to be unraveled. Like a tapestry woven with vibrant threads, these diverse texts unveil facets of human existence often veiled from our own unique vantage points. As we delve into stories stemming from divergent corners of the world, we begin to comprehend the enigmatic tapestry of human emotions, struggles, and triumphs. These literary works serve as portals enabling us to transcend our cultural boundaries, revealing the complexities and beauty embedded within the elusive mystery of human existence. With each novel, poem, or short story we dissect, we embarked on a journey of self-discovery and empathy, culminating in a profound appreciation for the symphony of voices that compose our collective humanity' + [char]11 + '' + [char]11 + 'In a world characterized by ever-evolving societal norms, the study of multicultural literature grants us the opportunity to engage in thoughtful examinations of the human condition across diverse cultural contexts. Through the characters we encounter and the landscapes they inhabit, we explore intricate relationships between identity, heritage, and personal narratives. By embarking on a literary odyssey through diverse perspectives, we cultivate an inclusive and nuanced understanding of the human spirit. These literary masterpieces encourage us to challenge prevailing viewpoints, question prevalent assumptions, and embrace a mindset open to the endless possibilities of human existence' + [char]11 + '' + [char]11 + 'Multicultural literature acts as a mirror reflecting the complexities and commonalities shared by all members of our global community. By fostering a deeper awareness of different cultures, we cultivate global citizens who actively seek understanding among diverse groups of people. As we analyze and interpret these literary offerings, we transcend mere appreciation and venture into actively dismantling barriers erected by prejudice, stereotypes, and narrow-mindedness. The profound lessons encapsulated within multicultural literature empower us to manifest a world united by empathy, tolerance, and an unwavering commitment to social justice.'
Replace-Text $bodyOld $bodyNew

# --- Summary paragraph ---
$summaryOld = 'Biodiversity, the immense diversity of life on Earth, encompasses the multitude of species inhabiting diverse ecosystems. Exploring biodiversity involves unraveling the intricate web of relationships among species and their environments. Understanding and preserving biodiversity is essential for maintaining ecological balance, supporting human well-being, and securing the future of life on Earth. Human activities pose significant threats to biodiversity, emphasizing the need for collective action to protect and restore ecosystems, ensuring the continued flourishing of life for generations to come.'
$summaryNew = 'Through captivating explorations of multicultural literature, we gain invaluable insights into the intricate threads that weave together the fabric of human experience. This genre unveils a panorama of diverse perspectives, fostering empathy, challenging societal norms, and cultivating a universal understanding of the human spirit. As we traverse these literary landscapes, we unearth the beauty of our collective existence and work towards building a more inclusive and tolerant world.'
Replace-Text $summaryOld $summaryNew

# --- Append a trailing empty paragraph right before the section break ---
$endOfStory = $d.Content.End
$tail = $d.Range($endOfStory - 1, $endOfStory - 1)
$tail.InsertParagraphAfter()

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
